# Sunay Nagpure Resume.docx edit:
#   1. Insert a new (empty-text) paragraph containing an explicit page
#      break right after the "Video editing" bullet and before the
#      "ACADEMIC PARTICIPATION" heading.
#   2. Word's re-layout after that insertion moves the recorded
#      "last rendered page break" marker: it now lands on the run that
#      starts the "ACADEMIC" heading (the first thing on the new page)
#      instead of on the "NPTEL Introduction to A..." run where it used
#      to sit.

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]10, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Step 1: insert the page-break paragraph right after "Video editing".
# ---------------------------------------------------------------------
$videoIdx = Get-ParaIndexByText $d "Video editing"
$videoRange = $d.Paragraphs.Item($videoIdx).Range
$videoRange.Collapse(0)

$pageBreakXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="61922507" w14:textId="6B929553" w:rsidR="006D0C02" w:rsidRPr="006D0C02" w:rsidRDefault="00F22653" w:rsidP="006D0C02"><w:pPr><w:pStyle w:val="Resume"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00A27CA1"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t>Video editing</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:br w:type="page"/></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$videoRange.InsertXML($pageBreakXml)

# ---------------------------------------------------------------------
# Step 2: add <w:lastRenderedPageBreak/> just before the "ACADEMIC" run
#         (first text on the newly created page).
# ---------------------------------------------------------------------
$academicIdx = Get-ParaIndexByText $d "ACADEMIC PARTICIPATION"
$academicRange = $d.Paragraphs.Item($academicIdx).Range
$academicRange.Collapse(0)

$academicXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="6BC6512E" w14:textId="3032C9F6" w:rsidR="00F81244" w:rsidRPr="00A27CA1" w:rsidRDefault="00F81244" w:rsidP="00F81244"><w:pPr><w:pStyle w:val="Resume"/><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="C9C9C9" w:themeFill="accent3" w:themeFillTint="99"/><w:spacing w:before="240" w:after="240"/><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00F81244"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>ACADEMIC</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00F81244"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>PARTICIPATION</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$academicRange.InsertXML($academicXml)

# ---------------------------------------------------------------------
# Step 3: remove the now-stale <w:lastRenderedPageBreak/> that used to
#         sit before the "NPTEL Introduction to A..." run.
# ---------------------------------------------------------------------
$nptelIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("NPTEL Introduction to Artificial Intelligence")) {
        $nptelIdx = $i
        break
    }
}
$nptelRange = $d.Paragraphs.Item($nptelIdx).Range
$nptelRange.Collapse(0)

$nptelXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="49DD0C67" w14:textId="41009076" w:rsidR="0067043F" w:rsidRDefault="0067043F" w:rsidP="00F81244"><w:pPr><w:pStyle w:val="Resume"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:spacing w:after="0"/><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="0067043F"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t>NPTEL Introduction to A</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t>rtificial Intelligence</w:t></w:r><w:r w:rsidRPr="0067043F"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Certificate</w:t></w:r><w:r w:rsidR="006D0C02"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> which had only 2.49% success rate</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$nptelRange.InsertXML($nptelXml)

Write-Output "Applied page-break insertion and lastRenderedPageBreak relocation."
